$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new log rows (167-207) capturing Admin activity on 2026-01-01
$logRows = @(
    @('2026-01-01 08:04:50', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2026-01-01 08:04:50', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2026-01-01 08:04:51', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2026-01-01 08:04:56', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:05:00', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:05:01', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:05:03', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:01', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2026-01-01 08:19:02', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2026-01-01 08:19:04', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:11', 'Admin', 'invoice', 'access_granted', 'Opened invoice page'),
    @('2026-01-01 08:19:31', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:42', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:50', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:53', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:56', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:56', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:19:58', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:20:00', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:07', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:07', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2026-01-01 08:27:08', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:13', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:13', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:16', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:27:35', 'Admin', 'invoice', 'access_granted', 'Opened invoice page'),
    @('2026-01-01 08:28:02', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:37:06', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:40:19', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:41:45', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2026-01-01 08:41:46', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2026-01-01 08:41:48', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:41:52', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:41:52', 'Admin', 'Quotation', 'quotation_created', 'Client: , Amount: 0.0'),
    @('2026-01-01 08:41:52', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:42:04', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:42:19', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:42:20', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:45:55', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:45:56', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2026-01-01 08:45:57', 'Admin', 'Quotation', 'quotation_created', 'Client: , Amount: 0.0')
)

$startRow = 167
for ($i = 0; $i -lt $logRows.Count; $i++) {
    $r = $startRow + $i
    $row = $logRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
